# Fruta / hortaliza, semanal
# Insert a new weekly record at row 176 (pushing the existing rows 176-200
# down to 177-201) and populate it with this week's Berenjena price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(176).Insert()

$ws.Range("A176").Value = 5
$ws.Range("B176").Value = "Macroferia Regional de Talca"
$ws.Range("C176").Value = "Maule"
$ws.Range("D176").Value = 45154
$ws.Range("E176").Value = 7
$ws.Range("F176").Value = 100112001
$ws.Range("G176").Value = "Berenjena"
$ws.Range("H176").Value = "Sin especificar"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 200
$ws.Range("K176").Value = 10000
$ws.Range("L176").Value = 10000
$ws.Range("M176").Value = 10000
$ws.Range("N176").Value = "$/caja 50 unidades"
$ws.Range("O176").Value = "Región de Arica y Parinacota"
$ws.Range("P176").Value = 200
$ws.Range("Q176").Value = 50
$ws.Range("R176").Value = "Hortaliza"
